# Turn off "Advance slide after N seconds" (removes the advTm="..." attribute
# from each slide's <p:transition>, both in the p14 mc:Choice and the legacy
# mc:Fallback) for every slide in the deck, leaving the rest of each slide's
# transition settings (e.g. p14:dur) untouched.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.SlideShowTransition.AdvanceOnTime = $false
}
